$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0.8430004128481028
    3  = 0.8149991627413542
    4  = 0.983795667930816
    5  = 0.5283687264620415
    6  = 0.9918647427602294
    7  = 0.9552030896374605
    8  = 0.8946756536048722
    9  = 0.183381994295238
    10 = 0.1999529876236299
    11 = 0.4282312392799456
    12 = 0.4251803619439751
    13 = 0.7413200363006528
    14 = 0.5718734448082547
    15 = 0.5894519237919532
    16 = 0.7562231448509459
    17 = 0.6999368509667938
    18 = 0.8003200282938729
    19 = 0.69545189377858
    20 = 0.5851580979618136
    21 = 0.8339375838625934
}

foreach ($row in $values.Keys) {
    $val = $values[$row]
    $ws.Range("C$row").Value = $val
    $ws.Range("D$row").Value = $val
}
